# Generate Report for Handoff
# Replaces the two "handed back" e2e sample files with two newly
# generated files that are now "Ready for handoff" (not yet handed back).
#
#   072e5562-1f19-4667-b4f4-49482e7ffe8d.md -> 71862875-3cce-4686-b8aa-acc814e9dcde.md
#   dbff55ea-36b5-4ee3-b886-a544221a78ad.md -> ffff47e1dfb2-2cbc-4270-b719-7b65a80e6f18.md

$wb = $excel.ActiveWorkbook

$ovw  = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$oldFile1 = "072e5562-1f19-4667-b4f4-49482e7ffe8d.md"
$oldFile2 = "dbff55ea-36b5-4ee3-b886-a544221a78ad.md"
$newFile1 = "71862875-3cce-4686-b8aa-acc814e9dcde.md"
$newFile2 = "ffff47e1dfb2-2cbc-4270-b719-7b65a80e6f18.md"

$oldXliffBase1 = "072e5562-1f19-4667-b4f4-49482e7ffe8d.f5285b9991bc0b3ad925509eb0ef97ddea9a3ef0"
$oldXliffBase2 = "dbff55ea-36b5-4ee3-b886-a544221a78ad.dee02c5f1fd82f6133af04c22b6692138f577c1b"
$newXliffBase  = "71862875-3cce-4686-b8aa-acc814e9dcde.6498da75e234a9e78b2721ed45546758b875fb12"

$newStatus = "Ready for handoff"

$newHoDateZhCn = "2016-08-15 22:58:24"
$newHoDateDeDe = "2016-08-15 22:58:29"
$zeroDate      = "0001-01-01 00:00:00"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acb16955c033ccde0bbec3579334f9e43841a156/e2e"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ovw.Range("A2").Value = $newFile1
$ovw.Range("B2").Value = "e2e\" + $newFile1
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("G2").Value = $newHoDateDeDe

$ovw.Range("A3").Value = $newFile2
$ovw.Range("B3").Value = "e2e\" + $newFile2
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus
$ovw.Range("G3").Value = $newHoDateDeDe

$ovw.Hyperlinks.Delete()
$ovw.Hyperlinks.Add($ovw.Range("B2"), $repoBase + "/" + $newFile1, "", "", "e2e\" + $newFile1)
$ovw.Hyperlinks.Add($ovw.Range("B3"), $repoBase + "/" + $newFile2, "", "", "e2e\" + $newFile2)

$ovw.Columns.Item(5).ColumnWidth = 16.25
$ovw.Columns.Item(6).ColumnWidth = 16.25

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn.Range("A2").Value = $newFile1
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("G2").Value = $newXliffBase + ".zh-cn.xlf"
$zhcn.Range("H2").Value = $newHoDateZhCn
$zhcn.Range("I2").Value = ""
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = $zeroDate

$zhcn.Range("A3").Value = $newFile2
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("F3").Value = "True"
$zhcn.Range("G3").Value = $newXliffBase + ".zh-cn.xlf"
$zhcn.Range("H3").Value = $newHoDateZhCn
$zhcn.Range("I3").Value = ""
$zhcn.Range("J3").Value = ""
$zhcn.Range("K3").Value = $zeroDate

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $repoBase + "/" + $newFile1, "", "", $newFile1)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $repoBase + "/" + $newFile2, "", "", $newFile2)

$zhcn.Columns.Item(3).ColumnWidth = 16.25
$zhcn.Columns.Item(9).ColumnWidth = 17.76
$zhcn.Columns.Item(10).ColumnWidth = 20.76

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede.Range("A2").Value = $newFile1
$dede.Range("C2").Value = $newStatus
$dede.Range("G2").Value = $newXliffBase + ".de-de.xlf"
$dede.Range("H2").Value = $newHoDateDeDe
$dede.Range("I2").Value = ""
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = $zeroDate

$dede.Range("A3").Value = $newFile2
$dede.Range("C3").Value = $newStatus
$dede.Range("F3").Value = "True"
$dede.Range("G3").Value = $newXliffBase + ".de-de.xlf"
$dede.Range("H3").Value = $newHoDateDeDe
$dede.Range("I3").Value = ""
$dede.Range("J3").Value = ""
$dede.Range("K3").Value = $zeroDate

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $repoBase + "/" + $newFile1, "", "", $newFile1)
$dede.Hyperlinks.Add($dede.Range("A3"), $repoBase + "/" + $newFile2, "", "", $newFile2)

$dede.Columns.Item(3).ColumnWidth = 16.25
$dede.Columns.Item(9).ColumnWidth = 17.76
$dede.Columns.Item(10).ColumnWidth = 20.76
